$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")

$ws.Cells.Item(420, 1).Value = "'2026-02-06"
$ws.Cells.Item(420, 2).Value = "10:14:51"
$ws.Cells.Item(420, 3).Value = "10:00"
$ws.Cells.Item(420, 4).Value = "Bathroom"
$ws.Cells.Item(420, 5).Value = "No Motion"
$ws.Cells.Item(420, 6).Value = "Inactive"

$ws.Cells.Item(421, 1).Value = "'2026-02-06"
$ws.Cells.Item(421, 2).Value = "10:14:52"
$ws.Cells.Item(421, 3).Value = "10:00"
$ws.Cells.Item(421, 4).Value = "Bathroom"
$ws.Cells.Item(421, 5).Value = "No Motion"
$ws.Cells.Item(421, 6).Value = "Inactive"

$ws.Cells.Item(422, 1).Value = "'2026-02-06"
$ws.Cells.Item(422, 2).Value = "10:14:56"
$ws.Cells.Item(422, 3).Value = "10:00"
$ws.Cells.Item(422, 4).Value = "Bathroom"
$ws.Cells.Item(422, 5).Value = "No Motion"
$ws.Cells.Item(422, 6).Value = "Inactive"

$ws.Cells.Item(423, 1).Value = "'2026-02-06"
$ws.Cells.Item(423, 2).Value = "10:15:01"
$ws.Cells.Item(423, 3).Value = "10:00"
$ws.Cells.Item(423, 4).Value = "Bathroom"
$ws.Cells.Item(423, 5).Value = "No Motion"
$ws.Cells.Item(423, 6).Value = "Inactive"

$ws.Cells.Item(424, 1).Value = "'2026-02-06"
$ws.Cells.Item(424, 2).Value = "10:15:06"
$ws.Cells.Item(424, 3).Value = "10:00"
$ws.Cells.Item(424, 4).Value = "Bathroom"
$ws.Cells.Item(424, 5).Value = "No Motion"
$ws.Cells.Item(424, 6).Value = "Inactive"

$ws.Cells.Item(425, 1).Value = "'2026-02-06"
$ws.Cells.Item(425, 2).Value = "10:15:11"
$ws.Cells.Item(425, 3).Value = "10:00"
$ws.Cells.Item(425, 4).Value = "Bathroom"
$ws.Cells.Item(425, 5).Value = "No Motion"
$ws.Cells.Item(425, 6).Value = "Inactive"

$ws.Cells.Item(426, 1).Value = "'2026-02-06"
$ws.Cells.Item(426, 2).Value = "10:15:14"
$ws.Cells.Item(426, 3).Value = "10:00"
$ws.Cells.Item(426, 4).Value = "Bathroom"
$ws.Cells.Item(426, 5).Value = "No Motion"
$ws.Cells.Item(426, 6).Value = "Inactive"

$ws.Cells.Item(427, 1).Value = "'2026-02-06"
$ws.Cells.Item(427, 2).Value = "10:15:19"
$ws.Cells.Item(427, 3).Value = "10:00"
$ws.Cells.Item(427, 4).Value = "Bathroom"
$ws.Cells.Item(427, 5).Value = "No Motion"
$ws.Cells.Item(427, 6).Value = "Inactive"

$ws.Cells.Item(428, 1).Value = "'2026-02-06"
$ws.Cells.Item(428, 2).Value = "10:15:24"
$ws.Cells.Item(428, 3).Value = "10:00"
$ws.Cells.Item(428, 4).Value = "Bathroom"
$ws.Cells.Item(428, 5).Value = "No Motion"
$ws.Cells.Item(428, 6).Value = "Inactive"

$ws.Cells.Item(429, 1).Value = "'2026-02-06"
$ws.Cells.Item(429, 2).Value = "10:15:29"
$ws.Cells.Item(429, 3).Value = "10:00"
$ws.Cells.Item(429, 4).Value = "Bathroom"
$ws.Cells.Item(429, 5).Value = "No Motion"
$ws.Cells.Item(429, 6).Value = "Inactive"

$ws.Cells.Item(430, 1).Value = "'2026-02-06"
$ws.Cells.Item(430, 2).Value = "10:15:36"
$ws.Cells.Item(430, 3).Value = "10:00"
$ws.Cells.Item(430, 4).Value = "Bathroom"
$ws.Cells.Item(430, 5).Value = "No Motion"
$ws.Cells.Item(430, 6).Value = "Inactive"

$ws.Cells.Item(431, 1).Value = "'2026-02-06"
$ws.Cells.Item(431, 2).Value = "10:15:41"
$ws.Cells.Item(431, 3).Value = "10:00"
$ws.Cells.Item(431, 4).Value = "Bathroom"
$ws.Cells.Item(431, 5).Value = "No Motion"
$ws.Cells.Item(431, 6).Value = "Inactive"

$ws.Cells.Item(432, 1).Value = "'2026-02-06"
$ws.Cells.Item(432, 2).Value = "10:15:46"
$ws.Cells.Item(432, 3).Value = "10:00"
$ws.Cells.Item(432, 4).Value = "Bathroom"
$ws.Cells.Item(432, 5).Value = "No Motion"
$ws.Cells.Item(432, 6).Value = "Inactive"

$ws.Cells.Item(433, 1).Value = "'2026-02-06"
$ws.Cells.Item(433, 2).Value = "10:15:49"
$ws.Cells.Item(433, 3).Value = "10:00"
$ws.Cells.Item(433, 4).Value = "Bathroom"
$ws.Cells.Item(433, 5).Value = "No Motion"
$ws.Cells.Item(433, 6).Value = "Inactive"

$ws = $wb.Worksheets.Item("Humidity")

$ws.Cells.Item(290, 1).Value = "'2026-02-06"
$ws.Cells.Item(290, 2).Value = "10:14:49"
$ws.Cells.Item(290, 3).Value = "10:00"
$ws.Cells.Item(290, 4).Value = "Bathroom"
$ws.Cells.Item(290, 5).Value = "'68.9%"
$ws.Cells.Item(290, 6).Value = "Active"

$ws.Cells.Item(291, 1).Value = "'2026-02-06"
$ws.Cells.Item(291, 2).Value = "10:14:54"
$ws.Cells.Item(291, 3).Value = "10:00"
$ws.Cells.Item(291, 4).Value = "Bathroom"
$ws.Cells.Item(291, 5).Value = "'69.0%"
$ws.Cells.Item(291, 6).Value = "Active"

$ws.Cells.Item(292, 1).Value = "'2026-02-06"
$ws.Cells.Item(292, 2).Value = "10:14:59"
$ws.Cells.Item(292, 3).Value = "10:00"
$ws.Cells.Item(292, 4).Value = "Bathroom"
$ws.Cells.Item(292, 5).Value = "'69.2%"
$ws.Cells.Item(292, 6).Value = "Active"

$ws.Cells.Item(293, 1).Value = "'2026-02-06"
$ws.Cells.Item(293, 2).Value = "10:15:04"
$ws.Cells.Item(293, 3).Value = "10:00"
$ws.Cells.Item(293, 4).Value = "Bathroom"
$ws.Cells.Item(293, 5).Value = "'69.1%"
$ws.Cells.Item(293, 6).Value = "Active"

$ws.Cells.Item(294, 1).Value = "'2026-02-06"
$ws.Cells.Item(294, 2).Value = "10:15:09"
$ws.Cells.Item(294, 3).Value = "10:00"
$ws.Cells.Item(294, 4).Value = "Bathroom"
$ws.Cells.Item(294, 5).Value = "'69.2%"
$ws.Cells.Item(294, 6).Value = "Active"

$ws.Cells.Item(295, 1).Value = "'2026-02-06"
$ws.Cells.Item(295, 2).Value = "10:15:34"
$ws.Cells.Item(295, 3).Value = "10:00"
$ws.Cells.Item(295, 4).Value = "Bathroom"
$ws.Cells.Item(295, 5).Value = "'69.0%"
$ws.Cells.Item(295, 6).Value = "Active"

$ws.Cells.Item(296, 1).Value = "'2026-02-06"
$ws.Cells.Item(296, 2).Value = "10:15:39"
$ws.Cells.Item(296, 3).Value = "10:00"
$ws.Cells.Item(296, 4).Value = "Bathroom"
$ws.Cells.Item(296, 5).Value = "'67.9%"
$ws.Cells.Item(296, 6).Value = "Active"

$ws.Cells.Item(297, 1).Value = "'2026-02-06"
$ws.Cells.Item(297, 2).Value = "10:15:44"
$ws.Cells.Item(297, 3).Value = "10:00"
$ws.Cells.Item(297, 4).Value = "Bathroom"
$ws.Cells.Item(297, 5).Value = "'68.9%"
$ws.Cells.Item(297, 6).Value = "Active"

$ws = $wb.Worksheets.Item("Temperature")

$ws.Cells.Item(290, 1).Value = "'2026-02-06"
$ws.Cells.Item(290, 2).Value = "10:14:50"
$ws.Cells.Item(290, 3).Value = "10:00"
$ws.Cells.Item(290, 4).Value = "Bathroom"
$ws.Cells.Item(290, 5).Value = "28.0C"
$ws.Cells.Item(290, 6).Value = "Active"

$ws.Cells.Item(291, 1).Value = "'2026-02-06"
$ws.Cells.Item(291, 2).Value = "10:14:55"
$ws.Cells.Item(291, 3).Value = "10:00"
$ws.Cells.Item(291, 4).Value = "Bathroom"
$ws.Cells.Item(291, 5).Value = "28.0C"
$ws.Cells.Item(291, 6).Value = "Active"

$ws.Cells.Item(292, 1).Value = "'2026-02-06"
$ws.Cells.Item(292, 2).Value = "10:15:00"
$ws.Cells.Item(292, 3).Value = "10:00"
$ws.Cells.Item(292, 4).Value = "Bathroom"
$ws.Cells.Item(292, 5).Value = "28.1C"
$ws.Cells.Item(292, 6).Value = "Active"

$ws.Cells.Item(293, 1).Value = "'2026-02-06"
$ws.Cells.Item(293, 2).Value = "10:15:05"
$ws.Cells.Item(293, 3).Value = "10:00"
$ws.Cells.Item(293, 4).Value = "Bathroom"
$ws.Cells.Item(293, 5).Value = "28.0C"
$ws.Cells.Item(293, 6).Value = "Active"

$ws.Cells.Item(294, 1).Value = "'2026-02-06"
$ws.Cells.Item(294, 2).Value = "10:15:10"
$ws.Cells.Item(294, 3).Value = "10:00"
$ws.Cells.Item(294, 4).Value = "Bathroom"
$ws.Cells.Item(294, 5).Value = "28.0C"
$ws.Cells.Item(294, 6).Value = "Active"

$ws.Cells.Item(295, 1).Value = "'2026-02-06"
$ws.Cells.Item(295, 2).Value = "10:15:35"
$ws.Cells.Item(295, 3).Value = "10:00"
$ws.Cells.Item(295, 4).Value = "Bathroom"
$ws.Cells.Item(295, 5).Value = "28.2C"
$ws.Cells.Item(295, 6).Value = "Active"

$ws.Cells.Item(296, 1).Value = "'2026-02-06"
$ws.Cells.Item(296, 2).Value = "10:15:40"
$ws.Cells.Item(296, 3).Value = "10:00"
$ws.Cells.Item(296, 4).Value = "Bathroom"
$ws.Cells.Item(296, 5).Value = "28.1C"
$ws.Cells.Item(296, 6).Value = "Active"

$ws.Cells.Item(297, 1).Value = "'2026-02-06"
$ws.Cells.Item(297, 2).Value = "10:15:45"
$ws.Cells.Item(297, 3).Value = "10:00"
$ws.Cells.Item(297, 4).Value = "Bathroom"
$ws.Cells.Item(297, 5).Value = "28.1C"
$ws.Cells.Item(297, 6).Value = "Active"

Write-Output "Added 14 PIR rows (420-433), 8 Humidity rows (290-297), 8 Temperature rows (290-297)"
